$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EffectiveDate (F) and PreviousExpDate (I) for both data rows
# from the old date string "07302023" to the new date string "08302023".
$ws.Range("F2").Value = "08302023"
$ws.Range("I2").Value = "08302023"
$ws.Range("F3").Value = "08302023"
$ws.Range("I3").Value = "08302023"

# Update the active cell selection to C8 (was G8).
$ws.Range("C8").Select()
